# Apply updated crypto price/volume values per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.237.44"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "'3.697.83"
$ws.Range("E3").Value = "  +7.89%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'582.67"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'178.01"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("D7").Value = "'3.691.67"
$ws.Range("E7").Value = "  +7.97%  "
$ws.Range("D8").Value = "'0.617"
$ws.Range("E8").Value = "  +4.05%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "'6.90"
$ws.Range("E11").Value = "  +27.51%  "
$ws.Range("E12").Value = "  +4.95%  "
$ws.Range("D13").Value = "'49.14"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "'0.0000287"
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("D15").Value = "'4.296.64"
$ws.Range("E15").Value = "  +8.17%  "
$ws.Range("D16").Value = "'679.00"
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("D17").Value = "'9.01"
$ws.Range("E17").Value = "  +4.65%  "
$ws.Range("D18").Value = "'3.698.08"
$ws.Range("E18").Value = "  +8.22%  "
$ws.Range("D19").Value = "'71.426.86"
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "'17.99"
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("E23").Value = "  +5.42%  "
$ws.Range("D24").Value = "'17.46"
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("D25").Value = "'102.20"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "'3.99"
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("D27").Value = "'2.84"
$ws.Range("E27").Value = "  +6.98%  "
$ws.Range("D28").Value = "'10.30"
$ws.Range("E28").Value = "  +7.80%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "'35.14"
$ws.Range("E30").Value = "  +4.95%  "
$ws.Range("D31").Value = "'3.42"
$ws.Range("E31").Value = "  +5.24%  "
$ws.Range("D32").Value = "'9.18"
$ws.Range("E32").Value = "  +4.93%  "
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").Value = "'7.57"
$ws.Range("E34").Value = "  +7.05%  "
$ws.Range("D35").Value = "'4.10"
$ws.Range("E35").Value = "  +10.67%  "
$ws.Range("D36").Value = "'582.21"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("E38").Value = "  +5.05%  "
$ws.Range("D39").Value = "'58.74"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").Value = "'3.677.24"
$ws.Range("E40").Value = "  +2.74%  "
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  +4.06%  "
$ws.Range("D43").Value = "'0.0459"
$ws.Range("E43").Value = "  +10.30%  "
$ws.Range("D44").Value = "'0.353"
$ws.Range("E44").Value = "  +6.42%  "
$ws.Range("D45").Value = "'0.0₃0771"
$ws.Range("E45").Value = "  +5.44%  "
$ws.Range("D46").Value = "'35.84"
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("D47").Value = "'2.78"
$ws.Range("E47").Value = "  +3.91%  "
$ws.Range("D48").Value = "'2.92"
$ws.Range("E48").Value = "  +10.49%  "
$ws.Range("D49").Value = "'0.134"
$ws.Range("E49").Value = "  +4.19%  "
$ws.Range("D50").Value = "'135.26"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("E51").Value = "  +10.60%  "
